$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.556.59"
$ws.Range("E2").Value = "  +1.39%  "
$ws.Range("D3").Value = "3.155.63"
$ws.Range("E3").Value = "  +0.86%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "594.84"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.80%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "147.14"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.65%  "
$ws.Range("D8").Value = "3.152.73"
$ws.Range("E8").Value = "  +1.09%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.533"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.42%  "
$ws.Range("E10").Value = "  -0.72%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.90"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +3.14%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.466"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.48%  "
$ws.Range("E13").Value = "  -1.96%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "37.49"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +3.53%  "
$ws.Range("D15").Value = "3.678.32"
$ws.Range("E15").Value = "  +0.82%  "
$ws.Range("E16").Value = "  -1.06%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "7.33"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +2.11%  "
$ws.Range("D18").Value = "64.322.11"
$ws.Range("E18").Value = "  +1.14%  "
$ws.Range("D19").Value = "3.161.43"
$ws.Range("E19").Value = "  +1.17%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "471.32"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.41%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.53"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.39%  "
$ws.Range("E22").Value = "  +0.48%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.60"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.62%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "13.19"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.11%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "81.80"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.72%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.32"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +5.25%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.00"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.02%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.61"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +6.91%  "
$ws.Range("B29").Value = "NEARProtocol"
$ws.Range("C29").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.49"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +8.13%  "
$ws.Range("B30").Value = "PancakeSwap"
$ws.Range("C30").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.74"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.95%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.26"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.22%  "
$ws.Range("E32").Value = "  +0.02%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "27.68"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.82%  "
$ws.Range("E34").Value = "  +0.95%  "
$ws.Range("D35").Value = "0.0₃0851"
$ws.Range("E35").Value = "  -0.52%  "
$ws.Range("E36").Value = "  +0.99%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.36"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.28%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.26"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.33%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.26"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.61%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "52.01"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.28%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "458.26"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.56%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "9.29"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +5.81%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.296"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +5.63%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0375"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.62%  "
$ws.Range("D45").Value = "2.945.19"
$ws.Range("E45").Value = "  +1.21%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "40.44"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +10.21%  "
$ws.Range("E47").Value = "  -1.23%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "129.42"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +3.09%  "
$ws.Range("E49").Value = "  -0.01%  "
$ws.Range("E50").Value = "  +3.06%  "
$ws.Range("E51").Value = "  -0.02%  "
